$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Get app crash logs in android phone" Heading4 bullet just
#    before the existing (empty) Heading4 paragraph that follows the
#    "... AppListAdapter" item.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*AppListAdapter*") {
        $target = $cand
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter() | Out-Null
    $newIndex = $target.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Range.Text = "Get app crash logs in android phone"
}

# ---------------------------------------------------------------------------
# 2) Move the <w:lastRenderedPageBreak/> marker from the "Dependency
#    Injection:" run up to the "References" run (the document grew a
#    paragraph above it, so the rendered page break now falls on the
#    "References" heading instead).
# ---------------------------------------------------------------------------
$refPara = $null
$depPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $t = $cand.Range.Text
    if ($t -eq "References`r") {
        $refPara = $cand
    } elseif ($t -eq "Dependency Injection:`r") {
        $depPara = $cand
    }
}

if (($refPara -ne $null) -and ($depPara -ne $null)) {
    $refXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2AF276F8" w14:textId="1D2A21E0" w:rsidR="00DD50D1" w:rsidRDefault="00DD50D1" w:rsidP="00DD50D1"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>References</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $refPara.Range.InsertXML($refXml) | Out-Null

    $depXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="54B57448" w14:textId="7ED3308C" w:rsidR="00DD50D1" w:rsidRDefault="00DD50D1" w:rsidP="00DD50D1"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Dependency Injection:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $depPara.Range.InsertXML($depXml) | Out-Null
}

Write-Output "edit complete"
